$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115, shifting existing rows 115-119 down to 116-120.
$ws.Range("A115:R115").EntireRow.Insert()

# Populate the newly inserted row 115 with the new weekly price record.
$ws.Range("A115").Value = 2
$ws.Range("B115").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C115").Value = "Coquimbo"
$ws.Range("D115").Value = 45021
$ws.Range("E115").Value = 4
$ws.Range("F115").Value = 100112030
$ws.Range("G115").Value = "Poroto granado"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 800
$ws.Range("K115").Value = 29000
$ws.Range("L115").Value = 30000
$ws.Range("M115").Value = 29500
$ws.Range("N115").Value = "`$/malla 25 kilos"
$ws.Range("O115").Value = "Provincia de Limarí"
$ws.Range("P115").Value = 1180
$ws.Range("Q115").Value = 25
$ws.Range("R115").Value = "Hortaliza"
